$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Energy")

# --- Labels / units -------------------------------------------------
# C6: clarify that the mobile-phone charging figure is a per-day value
$ws.Range("C6").Value = "Per day (full charge)"

# Fix the "Vaccum" typo for the Energy-sheet vacuum cleaner row, and
# keep it distinct from the (still misspelled) entry referenced by the
# Electricity sheet.
$ws.Range("A28").Value = "Vacuum cleaner (hoover)"

# Rows 25 & 26 (fridge / fridge-freezer) are fixed, always-on devices -
# add a "Per day" unit label next to them for clarity.
$ws.Range("C25").Value = "Per day"
$ws.Range("C26").Value = "Per day"

# --- Oven cycling factor --------------------------------------------
# Ovens don't run continuously at full power for the whole time they are
# on - apply a 0.55 cycling/duty-cycle factor to the energy estimate.
$ws.Range("E20").Formula = "=B20*D20*0.55"
$ws.Range("E20").NumberFormat = "0"

$ws.Range("E21").Formula = "=B21*D21*0.55"
$ws.Range("E21").NumberFormat = "0"

# --- View state -------------------------------------------------------
$ws.Activate()
$ws.Range("E20").Select()
